$wb = $excel.ActiveWorkbook

# Switch to the "Week 6" sheet (this becomes the active/selected tab)
$ws6 = $wb.Worksheets.Item("Week 6")
$ws6.Activate()

# Add the new "implementation" entry: 8 hours
$ws6.Range("A3").Value = "implementation"
$ws6.Range("B3").Value = 8

# Leave the cursor on the cell below the newly entered value
$ws6.Range("B4").Select()
